# Applies the "add MS" (Must Support) edit to the ROR PractitionerRole
# StructureDefinition spreadsheet:
#  - Flags a set of "Elements" rows as Must Support = "Y" (column H)
#  - Adds an AutoFilter over A1:AN79 with the two saved filter criteria
#  - Registers the (hidden) _FilterDatabase defined name that Excel writes
#    whenever a sheet carries an AutoFilter
#  - Hides all data rows (2-79) on the "Elements" sheet (Excel's IG
#    publisher workflow collapses the sheet to the filtered view) - this
#    must happen AFTER the AutoFilter call, which otherwise resets it
#  - Adds conditional formatting highlighting rows where Must Support isn't
#    "Y", and rows with a "Meaning When Missing" value
#  - Refreshes the recorded generation Date on the Metadata sheet

$wb = $excel.ActiveWorkbook
$elements = $wb.Worksheets.Item("Elements")
$metadata = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------
# 1. Mark the newly Must-Supported elements (column H = "Must Support?")
# ---------------------------------------------------------------------
$msRows = @(8, 13, 19, 20, 21, 26, 27, 28, 29, 31, 34, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 49, 50, 53, 54, 55, 57, 61, 64, 65, 66, 70, 71)
foreach ($r in $msRows) {
    $elements.Cells.Item($r, 8).Value = "Y"
}

# ---------------------------------------------------------------------
# 2. AutoFilter: column G (colId 6) excludes blank " ", column AA
#    (colId 26) keeps only blanks
# ---------------------------------------------------------------------
$fullRange = $elements.Range("A1:AN79")
$fullRange.AutoFilter(7, "<> ", 1)
$fullRange.AutoFilter(27, @(""))

# ---------------------------------------------------------------------
# 3. Hidden workbook-level _FilterDatabase name scoped to Elements,
#    matching what Excel persists alongside an AutoFilter
# ---------------------------------------------------------------------
$filterName = $elements.Names.Add("_xlnm._FilterDatabase", "=Elements!`$A`$1:`$AN`$79")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 4. Hide every element row (2-79), keeping the header row visible.
#    Must run after AutoFilter, which otherwise clears row visibility.
# ---------------------------------------------------------------------
$elements.Range("A2:A79").EntireRow.Hidden = $true

# ---------------------------------------------------------------------
# 5. Conditional formatting over the data rows: highlight rows missing
#    "Y" in Must Support, and italicize rows with a non-blank Meaning
#    When Missing value
# ---------------------------------------------------------------------
$cfRange = $elements.Range("A2:AI78")
$msRule = $cfRange.FormatConditions.Add(2, 0, '$G2<>"Y"')
$msRule.Interior.ColorIndex = 22

$mwmRule = $cfRange.FormatConditions.Add(2, 0, '$Q2<>""')
$mwmRule.Font.ColorIndex = 22
$mwmRule.Font.Italic = $true

# ---------------------------------------------------------------------
# 6. Refresh the recorded generation Date on the Metadata sheet
# ---------------------------------------------------------------------
$metadata.Cells.Item(8, 2).Value = "2024-02-07T14:28:15+00:00"
